$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 / column A: add the missing test-data value (keeps existing Hyperlink-style formatting, s="1")
$ws.Range("A3").Value = "testfsd34"

# New row 5: add the new regression test account e-mail, with hyperlink (mailto:) + hyperlink styling
$ws.Range("A5").Value = "testingdevbcregression1@yopmail.com"
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:testingdevbcregression1@yopmail.com")
$ws.Range("A5").Style = "Hyperlink"

# Resize columns to fit the new, wider content
$ws.Range("A1:F5").EntireColumn.AutoFit()

# Move the active selection the way it was left after the edit
$ws.Range("A6").Select()
